$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value into a B/C pair and make sure the cell keeps the
# same "column style" (s=2 for B, s=3 for C) that every other data row in
# this sheet uses. New cells created via COM sometimes inherit the wrong
# style from the <cols> definition (col A/B overlap), and any text that
# looks like a number/date gets auto-typed by Excel, silently dragging in
# a new (quotePrefix) style - so we always reassert the format afterwards
# by pasting it in from a known-good neighbour in the same column.
# ---------------------------------------------------------------------------
function Set-BC([string]$bAddr, [string]$cAddr, [string]$text, [string]$bFmtFrom, [string]$cFmtFrom) {
    $ws.Range($bAddr).Value2 = $text
    $ws.Range($bFmtFrom).Copy() | Out-Null
    $ws.Range($bAddr).PasteSpecial(-4122) | Out-Null

    $ws.Range($cAddr).Value2 = $text
    $ws.Range($cFmtFrom).Copy() | Out-Null
    $ws.Range($cAddr).PasteSpecial(-4122) | Out-Null
}

# 1) Name (English) lost the trailing "of materials"
Set-BC "B4" "C4" "Methods of electronic structure calculation" "B3" "C3"

# 2) Activation date changed. Note: B15/C15 already (erroneously) point at
#    the very same shared string as B8/C8 in the original workbook, so
#    updating both keeps that pre-existing aliasing behaviour intact.
#    Format is re-pasted in from row 7 (a clean, untouched row) because the
#    date-looking text would otherwise drag Excel's auto-typing/quote-prefix
#    style onto the very cell we're about to copy from.
Set-BC "B8" "C8" "'01/01/2023" "B7" "C7"
Set-BC "B15" "C15" "'01/01/2023" "B7" "C7"

# 3) Objectives (English) - new long-form paragraph added in B11/C11
$objectives = "Provide the student with a basic view of the main methods of theoretical determination of the electronic structure, focusing on crystalline solids, but also on molecules, two-dimensional materials and nanostructured materials. The main calculation method to be used in the course will be the Density Functional Theory (DFT), in some of its many variants. At the end of the course, the student will be able to determine material properties such as band structures, densities of states, elastic constants, and Fermi surfaces, using one or more of the methods and computer codes presented in class."
Set-BC "B11" "C11" $objectives "B10" "C10"

# 4) Short syllabus (English) - new cells B14/C14
$shortSyllabus = "Review of Quantum Mechanics; Review of Solid State Physics; Hartree-Fock Method; Density Functional Theory; Plane and pseudopotential wave methods; computer codes"
Set-BC "B14" "C14" $shortSyllabus "B13" "C13"

# 5) Full syllabus (English) - new cells B16/C16
$syllabus = "• Review of quantum mechanics: Schrödinger's equation; Hydrogen atom and atomic orbitals; Dirac notation; Variational principle; Linear combination of atomic orbitals. • Solid state physics review: Direct and reciprocal space; Bloch's Theorem; Brillouin zone; Energy bands and density of states; Fermi energy and Fermi surface; Free electrons Approximation. • Hartree-Fock method: Slater determinants; Hartree-Fock equation; Exchange and correlation potential; Self-consistent algorithm. • Density functional theory: Hohenberg-Kohn theorems; Kohn-Sham equations; Exchange and correlation functionals: LDA, GGA, etc. • Plane and pseudopotential wave methods: Plane wave bases; Pseudo-potentials; • Augmented and linearized plane wave bases: FP-LAPW method. • Computer codes: NWCHEM, Quantum Espresso, , Wien2k, exciting, VASP, etc."
Set-BC "B16" "C16" $syllabus "B13" "C13"
